$wb = $excel.ActiveWorkbook

# Rename identifiers to support the new multi-axle convention:
#   sAxleF -> sAxle1
#   Body_1Axle -> Body_Axle1
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A5").Value = "sAxle1"
    $ws.Range("H4").Value = "Body_Axle1"
}

# Restore each sheet's active-cell selection (bottom-right pane) and
# re-activate the sheets in the same order so the originally-selected tab
# stays selected.
$wsUnstable = $wb.Worksheets.Item("Trailer_Elula_Unstable")
$wsUnstable.Activate()
$wsUnstable.Range("H4").Select()

$wsThwala = $wb.Worksheets.Item("Trailer_Thwala")
$wsThwala.Activate()
$wsThwala.Range("H4").Select()

$wsElula = $wb.Worksheets.Item("Trailer_Elula")
$wsElula.Activate()
$wsElula.Range("M7").Select()
